# Add manual review notes (manualAudit / manualStatus) to the fastq tracking
# sheet, for the rows that were manually reviewed on 20200104.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that received a manual audit flag (column M) and a manual status
# note of "[512]" (column N).
$reviewedRows = @(9, 11, 16, 34, 45)

foreach ($r in $reviewedRows) {
    $ws.Cells.Item($r, 13).Value = 1          # column M -> manualAudit
    $ws.Cells.Item($r, 14).Value = "[512]"    # column N -> manualStatus
}

# Move the selection/view to where the last edit was made, matching the
# saved workbook view state.
$ws.Range("L45").Select()
